$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the "O bibliotecário pode adicionar..." paragraph (story summary)
# robustly via Find rather than a hard-coded paragraph index.
# ---------------------------------------------------------------------------
$finder = $d.Content
$found = $finder.Find.Execute(
    "O bibliotecário pode adicionar, excluir ou editar os gêneros de livros.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target paragraph."
}
$origIndex = $finder.Paragraphs.Item(1).Index

# ---------------------------------------------------------------------------
# Insert two new (initially blank) paragraphs immediately before it; each
# InsertParagraphBefore() call adds one blank paragraph right above the
# target while leaving the target paragraph itself untouched. (Re-fetch the
# paragraph range fresh from $d each time -- cached Range/Paragraph objects
# do not track subsequent document mutations.)
# ---------------------------------------------------------------------------
$d.Paragraphs.Item($origIndex).Range.InsertParagraphBefore() | Out-Null
$d.Paragraphs.Item($origIndex + 1).Range.InsertParagraphBefore() | Out-Null

$firstNewIndex = $origIndex
$secondNewIndex = $origIndex + 1
$storyIndex = $origIndex + 2

# ---------------------------------------------------------------------------
# Helper OOXML wrapper used to inject exact paragraph markup via
# Range.InsertXML (keeps formatting minimal / matches the authored edit).
# ---------------------------------------------------------------------------
$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function New-SimpleParagraphXml($text) {
    return '<w:p><w:pPr>' +
        '<w:spacing w:after="240" w:before="240" w:lineRule="auto"/>' +
        '<w:jc w:val="left"/>' +
        '<w:rPr/>' +
        '</w:pPr>' +
        '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr>' +
        '<w:t xml:space="preserve">' + $text + '</w:t></w:r></w:p>'
}

$text1 = "Na biblioteca, eu preciso manter os gêneros dos livros organizados. Mas às vezes surge um novo gênero e tenho que adicionar. Outras vezes, encontro erros e preciso editar. Também acontece de alguns gêneros não serem mais usados e eu preciso excluir."
$text2 = "O problema é garantir que a lista de gêneros esteja sempre atualizada para que os leitores encontrem os livros com facilidade."

$d.Paragraphs.Item($firstNewIndex).Range.InsertXML($xmlHeader + (New-SimpleParagraphXml $text1) + $xmlFooter)
$d.Paragraphs.Item($secondNewIndex).Range.InsertXML($xmlHeader + (New-SimpleParagraphXml $text2) + $xmlFooter)

# ---------------------------------------------------------------------------
# Simplify the original story paragraph: drop the descriptive run/text and
# collapse the paragraph-mark run properties down to an empty <w:rPr/>,
# keeping every other paragraph-formatting property unchanged.
# ---------------------------------------------------------------------------
$storyXml = '<w:p><w:pPr>' +
    '<w:keepNext w:val="0"/>' +
    '<w:keepLines w:val="0"/>' +
    '<w:pageBreakBefore w:val="0"/>' +
    '<w:widowControl w:val="1"/>' +
    '<w:pBdr>' +
    '<w:top w:space="0" w:sz="0" w:val="nil"/>' +
    '<w:left w:space="0" w:sz="0" w:val="nil"/>' +
    '<w:bottom w:space="0" w:sz="0" w:val="nil"/>' +
    '<w:right w:space="0" w:sz="0" w:val="nil"/>' +
    '<w:between w:space="0" w:sz="0" w:val="nil"/>' +
    '</w:pBdr>' +
    '<w:shd w:fill="auto" w:val="clear"/>' +
    '<w:spacing w:after="0" w:before="0" w:line="240" w:lineRule="auto"/>' +
    '<w:ind w:left="0" w:right="0" w:firstLine="0"/>' +
    '<w:jc w:val="left"/>' +
    '<w:rPr/>' +
    '</w:pPr>' +
    '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r>' +
    '</w:p>'

$d.Paragraphs.Item($storyIndex).Range.InsertXML($xmlHeader + $storyXml + $xmlFooter)

Write-Output "done"
